$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Range("E2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("C2").Value = [double]"0.02054416353671518"
$ws.Range("D2").Value = [double]"0.004102534773005262"
$ws.Range("E2").Value = [double]"0.006839779162894688"
$ws.Range("F2").Value = [double]"1.983307453579108E-05"
$ws.Range("G2").Value = [double]"3.78483761205306E-05"
$ws.Range("H2").Value = [double]"0.02182074568113457"
$ws.Range("I2").Value = [double]"0.01519077664026258"
$ws.Range("J2").Value = [double]"1.118703574221058E-09"
$ws.Range("B3").Value = [double]"0.02054416353671518"
$ws.Range("D3").Value = [double]"0.002463615943131936"
$ws.Range("E3").Value = [double]"0.002961048196503224"
$ws.Range("F3").Value = [double]"3.137397357222405E-05"
$ws.Range("G3").Value = [double]"9.95424009886392E-09"
$ws.Range("H3").Value = [double]"0.0003798629850140234"
$ws.Range("I3").Value = [double]"0.2452406342246212"
$ws.Range("J3").Value = [double]"1.224407686351014E-09"
$ws.Range("B4").Value = [double]"0.004102534773005262"
$ws.Range("C4").Value = [double]"0.002463615943131936"
$ws.Range("E4").Value = [double]"0.1130517233694399"
$ws.Range("F4").Value = [double]"0.2273095289132145"
$ws.Range("G4").Value = [double]"0.1208675045740597"
$ws.Range("H4").Value = [double]"0.6088297020432312"
$ws.Range("I4").Value = [double]"0.001446779506447671"
$ws.Range("J4").Value = [double]"1.341906807894588E-09"
$ws.Range("B5").Value = [double]"0.006839779162894688"
$ws.Range("C5").Value = [double]"0.002961048196503224"
$ws.Range("D5").Value = [double]"0.1130517233694399"
$ws.Range("F5").Value = [double]"0.4367912771902984"
$ws.Range("G5").Value = [double]"0.9539080177214705"
$ws.Range("H5").Value = [double]"0.6066196046388943"
$ws.Range("I5").Value = [double]"0.002042482977079763"
$ws.Range("J5").Value = [double]"1.733058403452503E-07"
$ws.Range("B6").Value = [double]"1.983307453579108E-05"
$ws.Range("C6").Value = [double]"3.137397357222405E-05"
$ws.Range("D6").Value = [double]"0.2273095289132145"
$ws.Range("E6").Value = [double]"0.4367912771902984"
$ws.Range("G6").Value = [double]"0.338115148263769"
$ws.Range("H6").Value = [double]"0.9398175514280553"
$ws.Range("I6").Value = [double]"1.37993871740516E-05"
$ws.Range("J6").Value = [double]"2.409578980788751E-08"
$ws.Range("B7").Value = [double]"3.78483761205306E-05"
$ws.Range("C7").Value = [double]"9.95424009886392E-09"
$ws.Range("D7").Value = [double]"0.1208675045740597"
$ws.Range("E7").Value = [double]"0.9539080177214705"
$ws.Range("F7").Value = [double]"0.338115148263769"
$ws.Range("H7").Value = [double]"0.2793491822940457"
$ws.Range("I7").Value = [double]"8.978023791073042E-08"
$ws.Range("J7").Value = [double]"2.873023352556459E-08"
$ws.Range("B8").Value = [double]"0.02182074568113457"
$ws.Range("C8").Value = [double]"0.0003798629850140234"
$ws.Range("D8").Value = [double]"0.6088297020432312"
$ws.Range("E8").Value = [double]"0.6066196046388943"
$ws.Range("F8").Value = [double]"0.9398175514280553"
$ws.Range("G8").Value = [double]"0.2793491822940457"
$ws.Range("I8").Value = [double]"0.000188894703086584"
$ws.Range("J8").Value = [double]"6.277102388985156E-07"
$ws.Range("B9").Value = [double]"0.01519077664026258"
$ws.Range("C9").Value = [double]"0.2452406342246212"
$ws.Range("D9").Value = [double]"0.001446779506447671"
$ws.Range("E9").Value = [double]"0.002042482977079763"
$ws.Range("F9").Value = [double]"1.37993871740516E-05"
$ws.Range("G9").Value = [double]"8.978023791073042E-08"
$ws.Range("H9").Value = [double]"0.000188894703086584"
$ws.Range("J9").Value = [double]"5.093762966978943E-09"
$ws.Range("B10").Value = [double]"1.118703574221058E-09"
$ws.Range("C10").Value = [double]"1.224407686351014E-09"
$ws.Range("D10").Value = [double]"1.341906807894588E-09"
$ws.Range("E10").Value = [double]"1.733058403452503E-07"
$ws.Range("F10").Value = [double]"2.409578980788751E-08"
$ws.Range("G10").Value = [double]"2.873023352556459E-08"
$ws.Range("H10").Value = [double]"6.277102388985156E-07"
$ws.Range("I10").Value = [double]"5.093762966978943E-09"

$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Range("C2").Value = [double]"2.610717141031786"
$ws.Range("D2").Value = [double]"-3.425139658064613"
$ws.Range("E2").Value = [double]"-3.168190191260131"
$ws.Range("F2").Value = [double]"-6.291540413517748"
$ws.Range("G2").Value = [double]"-5.912552836069583"
$ws.Range("H2").Value = [double]"-2.579732203178777"
$ws.Range("I2").Value = [double]"2.764999082356697"
$ws.Range("J2").Value = [double]"-14.13029932614793"
$ws.Range("B3").Value = [double]"-2.610717141031786"
$ws.Range("D3").Value = [double]"-3.681994804250974"
$ws.Range("E3").Value = [double]"-3.589237874196824"
$ws.Range("F3").Value = [double]"-6.021381501038252"
$ws.Range("G3").Value = [double]"-11.94190806746603"
$ws.Range("H3").Value = [double]"-4.643327564778173"
$ws.Range("I3").Value = [double]"1.212891665998232"
$ws.Range("J3").Value = [double]"-14.03349730670159"
$ws.Range("B4").Value = [double]"3.425139658064613"
$ws.Range("C4").Value = [double]"3.681994804250974"
$ws.Range("E4").Value = [double]"-1.690568841997594"
$ws.Range("F4").Value = [double]"-1.262742388823827"
$ws.Range("G4").Value = [double]"-1.651562663370604"
$ws.Range("H4").Value = [double]"-0.5234764216533269"
$ws.Range("I4").Value = [double]"3.951749147320574"
$ws.Range("J4").Value = [double]"-13.93584455663424"
$ws.Range("B5").Value = [double]"3.168190191260131"
$ws.Range("C5").Value = [double]"3.589237874196824"
$ws.Range("D5").Value = [double]"1.690568841997594"
$ws.Range("F5").Value = [double]"0.8005073618741826"
$ws.Range("G5").Value = [double]"0.05884403937409685"
$ws.Range("H5").Value = [double]"0.5267404159089561"
$ws.Range("I5").Value = [double]"3.77675076044754"
$ws.Range("J5").Value = [double]"-9.510493386921434"
$ws.Range("B6").Value = [double]"6.291540413517748"
$ws.Range("C6").Value = [double]"6.021381501038252"
$ws.Range("D6").Value = [double]"1.262742388823827"
$ws.Range("E6").Value = [double]"-0.8005073618741826"
$ws.Range("G6").Value = [double]"-0.9918280945160519"
$ws.Range("H6").Value = [double]"0.07686638543560266"
$ws.Range("I6").Value = [double]"6.509539756960056"
$ws.Range("J6").Value = [double]"-11.14188304278157"
$ws.Range("B7").Value = [double]"5.912552836069583"
$ws.Range("C7").Value = [double]"11.94190806746603"
$ws.Range("D7").Value = [double]"1.651562663370604"
$ws.Range("E7").Value = [double]"-0.05884403937409685"
$ws.Range("F7").Value = [double]"0.9918280945160519"
$ws.Range("H7").Value = [double]"1.125385968679645"
$ws.Range("I7").Value = [double]"10.03207067113413"
$ws.Range("J7").Value = [double]"-10.98803946835767"
$ws.Range("B8").Value = [double]"2.579732203178777"
$ws.Range("C8").Value = [double]"4.643327564778173"
$ws.Range("D8").Value = [double]"0.5234764216533269"
$ws.Range("E8").Value = [double]"-0.5267404159089561"
$ws.Range("F8").Value = [double]"-0.07686638543560266"
$ws.Range("G8").Value = [double]"-1.125385968679645"
$ws.Range("I8").Value = [double]"5.015819887476932"
$ws.Range("J8").Value = [double]"-8.54839354371728"
$ws.Range("B9").Value = [double]"-2.764999082356697"
$ws.Range("C9").Value = [double]"-1.212891665998232"
$ws.Range("D9").Value = [double]"-3.951749147320574"
$ws.Range("E9").Value = [double]"-3.77675076044754"
$ws.Range("F9").Value = [double]"-6.509539756960056"
$ws.Range("G9").Value = [double]"-10.03207067113413"
$ws.Range("H9").Value = [double]"-5.015819887476932"
$ws.Range("J9").Value = [double]"-12.57953253423835"
$ws.Range("B10").Value = [double]"14.13029932614793"
$ws.Range("C10").Value = [double]"14.03349730670159"
$ws.Range("D10").Value = [double]"13.93584455663424"
$ws.Range("E10").Value = [double]"9.510493386921434"
$ws.Range("F10").Value = [double]"11.14188304278157"
$ws.Range("G10").Value = [double]"10.98803946835767"
$ws.Range("H10").Value = [double]"8.54839354371728"
$ws.Range("I10").Value = [double]"12.57953253423835"

$ws = $wb.Worksheets.Item("Resumen")
$ws.Range("A2").Value = "Sieve Bootstrap"
$ws.Range("B2").Value = 4
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = [double]"0.6444787228878719"
$ws.Range("A3").Value = "DeepAR"
$ws.Range("B3").Value = 4
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = [double]"0.6332213953123457"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = [double]"37.5"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 7
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 7
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 5

